$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5856.2
$ws.Range("I64").Value = 5820.25
$ws.Range("J64").Value = 6000
$ws.Range("K64").Value = 5820.25
$ws.Range("L64").Value = 6000
$ws.Range("M64").Value = -5572.25
$ws.Range("N64").Value = -6496
$ws.Range("H67").Value = 5856.2
$ws.Range("I67").Value = 5820.25
$ws.Range("J67").Value = 6000
$ws.Range("K67").Value = 5820.25
$ws.Range("L67").Value = 6000
$ws.Range("M67").Value = -4962.25
$ws.Range("N67").Value = -7716
$ws.Range("H86").Value = 6711.909
$ws.Range("I86").Value = 7453.4443
$ws.Range("K86").Value = 7453.4443
$ws.Range("M86").Value = -6330.4443
$ws.Range("H89").Value = 6711.909
$ws.Range("I89").Value = 7453.4443
$ws.Range("K89").Value = 37267.2215
$ws.Range("M89").Value = -31651.2215
$ws.Range("H95").Value = 189950
$ws.Range("J95").Value = 189950
$ws.Range("L95").Value = 189950
$ws.Range("N95").Value = -195442
$ws.Range("H106").Value = 1600.7142
$ws.Range("I106").Value = 941
$ws.Range("K106").Value = 941
$ws.Range("M106").Value = -310
$ws.Range("H113").Value = 7697.9
$ws.Range("J113").Value = 4506.6665
$ws.Range("L113").Value = 4506.6665
$ws.Range("N113").Value = -11014.6665
$ws.Range("H135").Value = 2299.4
$ws.Range("I135").Value = 2374.25
$ws.Range("K135").Value = 21368.25
$ws.Range("M135").Value = -18833.25
$ws.Range("H137").Value = 7952.159
$ws.Range("I137").Value = 3776
$ws.Range("K137").Value = 11328
$ws.Range("M137").Value = -8778
$ws.Range("H138").Value = 3760.8438
$ws.Range("J138").Value = 3881.611
$ws.Range("L138").Value = 11644.833
$ws.Range("N138").Value = -21924.833

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3617.2576
$ws.Range("I32").Value = 1356.102
$ws.Range("K32").Value = 1356.102
$ws.Range("M32").Value = -1069.102
$ws.Range("H74").Value = 8174.3145
$ws.Range("I74").Value = 2582.9333
$ws.Range("K74").Value = 2582.9333
$ws.Range("M74").Value = -1708.9333
$ws.Range("H77").Value = 8174.3145
$ws.Range("I77").Value = 2582.9333
$ws.Range("K77").Value = 12914.6665
$ws.Range("M77").Value = -8546.666500000001
$ws.Range("H95").Value = 57403.25
$ws.Range("J95").Value = 57403.25
$ws.Range("L95").Value = 57403.25
$ws.Range("N95").Value = -62895.25
$ws.Range("H97").Value = 3671.5557
$ws.Range("I97").Value = 3671.5557
$ws.Range("K97").Value = 3671.5557
$ws.Range("M97").Value = -3175.5557
$ws.Range("H122").Value = 1431638.9
$ws.Range("I122").Value = 1668512
$ws.Range("J122").Value = 10400
$ws.Range("K122").Value = 5005536
$ws.Range("L122").Value = 31200
$ws.Range("M122").Value = -5003086
$ws.Range("N122").Value = -36100
$ws.Range("H133").Value = 90574.11
$ws.Range("J133").Value = 90574.11
$ws.Range("L133").Value = 90574.11
$ws.Range("N133").Value = -95634.11

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1667.5
$ws.Range("I22").Value = 1901
$ws.Range("K22").Value = 1901
$ws.Range("M22").Value = -1728
$ws.Range("H80").Value = 929.7619
$ws.Range("I80").Value = 920.44446
$ws.Range("J80").Value = 936.75
$ws.Range("K80").Value = 920.44446
$ws.Range("L80").Value = 936.75
$ws.Range("M80").Value = 77.55553999999995
$ws.Range("N80").Value = -2932.75
$ws.Range("H83").Value = 929.7619
$ws.Range("I83").Value = 920.44446
$ws.Range("J83").Value = 936.75
$ws.Range("K83").Value = 4602.2223
$ws.Range("L83").Value = 4683.75
$ws.Range("M83").Value = 389.7776999999996
$ws.Range("N83").Value = -14667.75
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 129389.37
$ws.Range("I31").Value = 224841.33
$ws.Range("K31").Value = 224841.33
$ws.Range("M31").Value = -224546.33
$ws.Range("H34").Value = 129389.37
$ws.Range("I34").Value = 224841.33
$ws.Range("K34").Value = 224841.33
$ws.Range("M34").Value = -224639.33
$ws.Range("H58").Value = 9119.518
$ws.Range("I58").Value = 4773.1
$ws.Range("K58").Value = 4773.1
$ws.Range("M58").Value = -4570.1
$ws.Range("H62").Value = 2994
$ws.Range("I62").Value = 2994
$ws.Range("J62").Value = 2994
$ws.Range("K62").Value = 2994
$ws.Range("L62").Value = 2994
$ws.Range("M62").Value = -2370
$ws.Range("N62").Value = -4242
$ws.Range("H65").Value = 2994
$ws.Range("I65").Value = 2994
$ws.Range("J65").Value = 2994
$ws.Range("K65").Value = 14970
$ws.Range("L65").Value = 14970
$ws.Range("M65").Value = -11850
$ws.Range("N65").Value = -21210
$ws.Range("H94").Value = 11402
$ws.Range("J94").Value = 12234.223
$ws.Range("L94").Value = 12234.223
$ws.Range("N94").Value = -13136.223
$ws.Range("H132").Value = 1644677
$ws.Range("I132").Value = 3257.25
$ws.Range("K132").Value = 9771.75
$ws.Range("M132").Value = -7241.75
$ws.Range("H136").Value = 9119.518
$ws.Range("I136").Value = 4773.1
$ws.Range("K136").Value = 14319.3
$ws.Range("M136").Value = -11769.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 482.92856
$ws.Range("I5").Value = 476.375
$ws.Range("J5").Value = 491.66666
$ws.Range("K5").Value = 1429.125
$ws.Range("L5").Value = 1474.99998
$ws.Range("M5").Value = -1317.125
$ws.Range("N5").Value = -1698.99998
$ws.Range("H38").Value = 49.642857
$ws.Range("J38").Value = 70.77778000000001
$ws.Range("L38").Value = 212.33334
$ws.Range("N38").Value = -906.33334
$ws.Range("H122").Value = 16554634
$ws.Range("J122").Value = 7084457
$ws.Range("L122").Value = 63760113
$ws.Range("N122").Value = -63765013
$ws.Range("H129").Value = 2346.4614
$ws.Range("I129").Value = 838
$ws.Range("J129").Value = 3289.25
$ws.Range("K129").Value = 2514
$ws.Range("L129").Value = 9867.75
$ws.Range("M129").Value = 2486
$ws.Range("N129").Value = -19867.75
$ws.Range("H131").Value = 1461.89
$ws.Range("I131").Value = 1025.8
$ws.Range("J131").Value = 1484.8422
$ws.Range("K131").Value = 3077.4
$ws.Range("L131").Value = 4454.5266
$ws.Range("M131").Value = 1962.6
$ws.Range("N131").Value = -14534.5266
$ws.Range("H135").Value = 482.92856
$ws.Range("I135").Value = 476.375
$ws.Range("J135").Value = 491.66666
$ws.Range("K135").Value = 4287.375
$ws.Range("L135").Value = 4424.99994
$ws.Range("M135").Value = -1752.375
$ws.Range("N135").Value = -9494.99994
$ws.Range("H140").Value = 2166.2856
$ws.Range("I140").Value = 1602.5454
$ws.Range("J140").Value = 4233.3335
$ws.Range("K140").Value = 4807.6362
$ws.Range("L140").Value = 12700.0005
$ws.Range("M140").Value = 372.3638000000001
$ws.Range("N140").Value = -23060.0005
$ws.Range("H141").Value = 5798.9
$ws.Range("I141").Value = 4748.625
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 14245.875
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -9065.875
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 80000
$ws.Range("J52").Value = 80000
$ws.Range("L52").Value = 80000
$ws.Range("N52").Value = -80518
$ws.Range("H80").Value = 2004.75
$ws.Range("I80").Value = 1826.8572
$ws.Range("K80").Value = 1826.8572
$ws.Range("M80").Value = -828.8571999999999
$ws.Range("H83").Value = 2004.75
$ws.Range("I83").Value = 1826.8572
$ws.Range("K83").Value = 9134.286
$ws.Range("M83").Value = -4142.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7204.25
$ws.Range("I7").Value = 7842.5713
$ws.Range("K7").Value = 7842.5713
$ws.Range("M7").Value = -7730.5713
$ws.Range("H46").Value = 3815.3684
$ws.Range("J46").Value = 3896.1333
$ws.Range("L46").Value = 3896.1333
$ws.Range("N46").Value = -4272.1333
$ws.Range("H122").Value = 5415.8887
$ws.Range("I122").Value = 6055.7856
$ws.Range("K122").Value = 18167.3568
$ws.Range("M122").Value = -15717.3568
$ws.Range("H126").Value = 7204.25
$ws.Range("I126").Value = 7842.5713
$ws.Range("K126").Value = 23527.7139
$ws.Range("M126").Value = -21057.7139
$ws.Range("H136").Value = 1233571.6
$ws.Range("I136").Value = 15304.333
$ws.Range("K136").Value = 45912.999
$ws.Range("M136").Value = -43362.999
$ws.Range("H137").Value = 97375
$ws.Range("J137").Value = 97375
$ws.Range("L137").Value = 97375
$ws.Range("N137").Value = -107575

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 16184.926
$ws.Range("I126").Value = 19795.143
$ws.Range("J126").Value = 3549.1667
$ws.Range("K126").Value = 59385.429
$ws.Range("L126").Value = 10647.5001
$ws.Range("M126").Value = -56915.429
$ws.Range("N126").Value = -15587.5001
